$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "564.55", "7.00").
# Force the cell format to Text first so Excel does not silently convert
# them to numbers (which would drop trailing zeros / change precision).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.319.09"
$ws.Range("E2").Value = "  +4.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.997.27"
$ws.Range("E3").Value = "  +3.92%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.55"
$ws.Range("E5").Value = "  +3.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.20"
$ws.Range("E6").Value = "  +13.05%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.520"
$ws.Range("E8").Value = "  +4.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.991.11"
$ws.Range("E9").Value = "  +3.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.133"
$ws.Range("E10").Value = "  +10.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.04"
$ws.Range("E11").Value = "  +9.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  +4.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  +10.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.83"
$ws.Range("E14").Value = "  +4.29%  "
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.496.30"
$ws.Range("E16").Value = "  +4.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.00"
$ws.Range("E17").Value = "  +7.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.996.35"
$ws.Range("E18").Value = "  +3.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.345.34"
$ws.Range("E19").Value = "  +4.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "429.95"
$ws.Range("E20").Value = "  +6.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("E21").Value = "  +6.70%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.712"
$ws.Range("E22").Value = "  +6.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.18"
$ws.Range("E23").Value = "  +5.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.49"
$ws.Range("E24").Value = "  +6.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.78"
$ws.Range("E25").Value = "  +4.63%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.13"
$ws.Range("E28").Value = "  +11.36%  "
$ws.Range("E29").Value = "  +4.79%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.80"
$ws.Range("E30").Value = "  +8.81%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.72"
$ws.Range("E31").Value = "  +4.56%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.16"
$ws.Range("E32").Value = "  +3.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0995"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("E34").Value = "  +10.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0₃0768"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.78"
$ws.Range("E36").Value = "  +7.58%  "
$ws.Range("E37").Value = "  +4.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.17"
$ws.Range("E38").Value = "  +2.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.67"
$ws.Range("E39").Value = "  +4.73%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.75"
$ws.Range("E40").Value = "  +15.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "406.54"
$ws.Range("E41").Value = "  +13.85%  "
$ws.Range("E42").Value = "  +4.77%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.757.27"
$ws.Range("E43").Value = "  +5.07%  "
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.252"
$ws.Range("E45").Value = "  +10.69%  "
$ws.Range("E46").Value = "  +0.04%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "123.92"
$ws.Range("E47").Value = "  +3.64%  "
$ws.Range("E48").Value = "  +4.82%  "
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.49"
$ws.Range("E50").Value = "  +20.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.57"
$ws.Range("E51").Value = "  +3.84%  "
